# Forest data update - 2025-11-26 12:20
# ------------------------------------------------------------------
# "Previously added" sheet keeps a running history; "New" sheet always
# holds the 3 most-recently scraped listings. On every scrape run the
# current contents of "New" get appended to the bottom of
# "Previously added", and "New" is overwritten with the 3 fresh
# listings.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# xlPasteSpecial constants used below
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --------------------------------------------------------------
# Helper: write a cadastre-style value that LOOKS numeric but must
# stay a text cell (cadastre numbers, not numbers). We build it with
# a throw-away formula that evaluates to a string, then paste-special
# "Values" back onto itself: this converts the formula to a plain
# shared-string cell without ever touching NumberFormat (which would
# otherwise create a brand-new, permanently-unused cell style).
# --------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

# ================================================================
# STEP 1 - move the CURRENT "New" rows (2-4) onto the end of
# "Previously added" (rows 281-283), verbatim, BEFORE we touch "New".
# Doing this first means the shared strings already used by these
# three rows stay alive (now referenced twice) instead of being
# overwritten in step 2.
# ================================================================

$destRow = 281

# --- old New!row2 -> Previously added!row281 ---
Set-TextValue $ws1.Range("E$destRow") "38460040026"
$ws1.Range("A$destRow").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/balvu-pag/akimk.html"
$ws1.Range("B$destRow").Value = "15 000 €"
$ws1.Range("C$destRow").Value = "Balvi un raj."
$ws1.Range("D$destRow").Value = "1 ha."
$ws1.Range("F$destRow").Value = 45985.790972222225

# --- old New!row3 -> Previously added!row282 ---
$destRow = 282
Set-TextValue $ws1.Range("E$destRow") "46600010112"
$ws1.Range("A$destRow").Value = "https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/dobeles-pag/lfnll.html"
$ws1.Range("B$destRow").Value = "35 000 €"
$ws1.Range("C$destRow").Value = "Dobele un raj."
$ws1.Range("D$destRow").Value = "6 ha."
$ws1.Range("F$destRow").Value = 45985.72777777778

# --- old New!row4 -> Previously added!row283 ---
$destRow = 283
$ws1.Range("A$destRow").Value = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/aronas-pag/cghlhb.html"
$ws1.Range("B$destRow").Value = "80 000 €"
$ws1.Range("C$destRow").Value = "Madona un raj."
$ws1.Range("D$destRow").Value = "137000 m²"
$ws1.Range("E$destRow").Value = ""
$ws1.Range("F$destRow").Value = 45985.65902777778

# Re-apply the formatting of the last pre-existing data row (280) onto
# the 3 new rows, so the new cells pick up the same styles (s="3" for
# the link column, s="4" for plain text columns, s="2" for the date
# column) instead of whatever ad-hoc style the Value assignments above
# left behind.
$ws1.Range("A280:F280").Copy()
$ws1.Range("A281:F283").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Hyperlink the link column (A) of the 3 new rows.
$ws1.Hyperlinks.Add($ws1.Range("A281"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/balvu-pag/akimk.html")
$ws1.Hyperlinks.Add($ws1.Range("A282"), "https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/dobeles-pag/lfnll.html")
$ws1.Hyperlinks.Add($ws1.Range("A283"), "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/aronas-pag/cghlhb.html")

# Hyperlinks.Add forces its own built-in "Hyperlink" style onto the
# cell; restore the correct style (s="3") by re-pasting the row-280
# formats once more.
$ws1.Range("A280:F280").Copy()
$ws1.Range("A281:F283").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ================================================================
# STEP 2 - overwrite "New" (rows 2-4) with the 3 freshly scraped
# listings.
# ================================================================

# --- New!row2 ---
$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/susaju-pag/cixen.html"
$ws2.Range("B2").Value = "6 000 €"
$ws2.Range("C2").Value = "Balvi un raj."
$ws2.Range("D2").Value = "1 ha."
$ws2.Range("E2").Value = ""
$ws2.Range("F2").Value = 45986.88055555556

# --- New!row3 ---
$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/davinu-pag/emmkl.html"
$ws2.Range("B3").Value = "21 000 €"
$ws2.Range("C3").Value = "Bauska un raj."
$ws2.Range("D3").Value = "7 ha."
Set-TextValue $ws2.Range("E3") "40560070003"
$ws2.Range("F3").Value = 45987.43611111111

# --- New!row4 ---
$ws2.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/pusas-pag/okndh.html"
$ws2.Range("B4").Value = "4 400 €"
$ws2.Range("C4").Value = "Rēzekne un raj."
$ws2.Range("D4").Value = "1 ha."
Set-TextValue $ws2.Range("E4") "78800040225"
$ws2.Range("F4").Value = 45986.683333333334

# Restore styling on the 3 rows we just rewrote (same reasoning as
# above for sheet 1).
$ws2.Range("A2:F2").Copy()
$ws2.Range("A2:F4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Replace the 3 old hyperlinks on "New" with the 3 new ones.
$ws2.Range("A2:A4").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/susaju-pag/cixen.html")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/davinu-pag/emmkl.html")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/pusas-pag/okndh.html")

# Restore the correct style once more after Hyperlinks.Add.
$ws2.Range("A2:F2").Copy()
$ws2.Range("A2:F4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Write-Host "Forest data updated."
